$d = $word.ActiveDocument

$replacements = @(
    @("501×4=", "105×2="),
    @("905×9=", "295×4="),
    @("735×9=", "133×7="),
    @("534×2=", "803×6="),
    @("869×9=", "383×5="),
    @("350×8=", "644×7="),
    @("115×7=", "599×4="),
    @("772×6=", "298×5="),
    @("429×5=", "400×3="),
    @("674×5=", "646×7="),
    @("399×8=", "635×4="),
    @("977×2=", "543×3="),
    @("229×5=", "496×6="),
    @("806×8=", "425×6="),
    @("167×5=", "704×5="),
    @("364×2=", "831×3="),
    @("303×9=", "157×8="),
    @("708×2=", "151×2="),
    @("972×4=", "920×4="),
    @("236×9=", "295×9="),
    @("891×6=", "246×7="),
    @("898×8=", "999×7="),
    @("354×9=", "595×4="),
    @("109×5=", "987×2="),
    @("238×7=", "495×5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
